$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 3242.4285  # H9
$ws.Cells.Item(9, 9).Value = 252.5  # I9
$ws.Cells.Item(9, 11).Value = 252.5  # K9
$ws.Cells.Item(9, 13).Value = -83.5  # M9

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 597.75  # H28
$ws.Cells.Item(28, 9).Value = 647.0714  # I28
$ws.Cells.Item(28, 11).Value = 647.0714  # K28
$ws.Cells.Item(28, 13).Value = -162.0714  # M28

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 2512.8667  # H38
$ws.Cells.Item(38, 9).Value = 790.4545000000001  # I38
$ws.Cells.Item(38, 11).Value = 2371.3635  # K38
$ws.Cells.Item(38, 13).Value = -1999.3635  # M38

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 20166.666  # H69
$ws.Cells.Item(69, 9).Value = 31000  # I69
$ws.Cells.Item(69, 11).Value = 93000  # K69
$ws.Cells.Item(69, 13).Value = -92126  # M69

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 3165.6667  # H70
$ws.Cells.Item(70, 10).Value = 3165.6667  # J70
$ws.Cells.Item(70, 12).Value = 9497.000100000001  # L70
$ws.Cells.Item(70, 14).Value = -10037.0001  # N70

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 20166.666  # H72
$ws.Cells.Item(72, 9).Value = 31000  # I72
$ws.Cells.Item(72, 11).Value = 279000  # K72
$ws.Cells.Item(72, 13).Value = -274632  # M72

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 3165.6667  # H73
$ws.Cells.Item(73, 10).Value = 3165.6667  # J73
$ws.Cells.Item(73, 12).Value = 9497.000100000001  # L73
$ws.Cells.Item(73, 14).Value = -11369.0001  # N73

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 1082.2222  # H94
$ws.Cells.Item(94, 9).Value = 1082.2222  # I94
$ws.Cells.Item(94, 11).Value = 1082.2222  # K94
$ws.Cells.Item(94, 13).Value = -631.2221999999999  # M94

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 3986.7778  # H135
$ws.Cells.Item(135, 9).Value = 4046.8823  # I135
$ws.Cells.Item(135, 11).Value = 36421.9407  # K135
$ws.Cells.Item(135, 13).Value = -33886.9407  # M135

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 1725.909  # H141
$ws.Cells.Item(141, 9).Value = 1725.909  # I141
$ws.Cells.Item(141, 10).Value = 0  # J141
$ws.Cells.Item(141, 11).Value = 5177.727000000001  # K141
$ws.Cells.Item(141, 12).Value = 0  # L141
$ws.Cells.Item(141, 13).Value = 2.272999999999229  # M141
$ws.Cells.Item(141, 14).ClearContents()  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 11248.263  # H63
$ws.Cells.Item(63, 9).Value = 4113.727  # I63
$ws.Cells.Item(63, 11).Value = 4113.727  # K63
$ws.Cells.Item(63, 13).Value = -3427.727  # M63

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 11248.263  # H66
$ws.Cells.Item(66, 9).Value = 4113.727  # I66
$ws.Cells.Item(66, 11).Value = 20568.635  # K66
$ws.Cells.Item(66, 13).Value = -17136.635  # M66

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2105.1667  # H122
$ws.Cells.Item(122, 9).Value = 1842.409  # I122
$ws.Cells.Item(122, 11).Value = 5527.227000000001  # K122
$ws.Cells.Item(122, 13).Value = -3077.227000000001  # M122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 19866.334  # H2
$ws.Cells.Item(2, 10).Value = 4799.5  # J2
$ws.Cells.Item(2, 12).Value = 4799.5  # L2
$ws.Cells.Item(2, 14).Value = -5025.5  # N2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 9439.031000000001  # H107
$ws.Cells.Item(107, 9).Value = 11890.131  # I107
$ws.Cells.Item(107, 10).Value = 3175.111  # J107
$ws.Cells.Item(107, 11).Value = 11890.131  # K107
$ws.Cells.Item(107, 12).Value = 3175.111  # L107
$ws.Cells.Item(107, 13).Value = -9970.130999999999  # M107
$ws.Cells.Item(107, 14).Value = -7015.111  # N107

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 99996  # H140
$ws.Cells.Item(140, 10).Value = 99996  # J140
$ws.Cells.Item(140, 12).Value = 99996  # L140
$ws.Cells.Item(140, 14).Value = -110356  # N140

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2882.2307  # H16
$ws.Cells.Item(16, 10).Value = 5018  # J16
$ws.Cells.Item(16, 12).Value = 5018  # L16
$ws.Cells.Item(16, 14).Value = -5592  # N16

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 50000  # H80
$ws.Cells.Item(80, 9).Value = 50000  # I80
$ws.Cells.Item(80, 11).Value = 50000  # K80
$ws.Cells.Item(80, 13).Value = -48877  # M80

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(83, 8).Value = 50000  # H83
$ws.Cells.Item(83, 9).Value = 50000  # I83
$ws.Cells.Item(83, 11).Value = 150000  # K83
$ws.Cells.Item(83, 13).Value = -144384  # M83

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1147.3334  # H94
$ws.Cells.Item(94, 10).Value = 1493.1111  # J94
$ws.Cells.Item(94, 12).Value = 1493.1111  # L94
$ws.Cells.Item(94, 14).Value = -2395.1111  # N94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 74000  # H97
$ws.Cells.Item(97, 10).Value = 74000  # J97
$ws.Cells.Item(97, 12).Value = 74000  # L97
$ws.Cells.Item(97, 14).Value = -75982  # N97

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(102, 8).Value = 82000  # H102
$ws.Cells.Item(102, 10).Value = 82000  # J102
$ws.Cells.Item(102, 12).Value = 82000  # L102
$ws.Cells.Item(102, 14).Value = -86868  # N102

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 2882.2307  # H113
$ws.Cells.Item(113, 10).Value = 5018  # J113
$ws.Cells.Item(113, 12).Value = 5018  # L113
$ws.Cells.Item(113, 14).Value = -9358  # N113

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 8467532  # H4
$ws.Cells.Item(4, 9).Value = 12286284  # I4
$ws.Cells.Item(4, 11).Value = 36858852  # K4
$ws.Cells.Item(4, 13).Value = -36858740  # M4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 1003451.9  # H103
$ws.Cells.Item(103, 9).Value = 2000264.6  # I103
$ws.Cells.Item(103, 11).Value = 6000793.800000001  # K103
$ws.Cells.Item(103, 13).Value = -5999914.800000001  # M103

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 7074035.5  # H131
$ws.Cells.Item(131, 10).Value = 148210.72  # J131
$ws.Cells.Item(131, 12).Value = 444632.16  # L131
$ws.Cells.Item(131, 14).Value = -454712.16  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 27333.334  # H98
$ws.Cells.Item(98, 10).Value = 27333.334  # J98
$ws.Cells.Item(98, 12).Value = 27333.334  # L98
$ws.Cells.Item(98, 14).Value = -33323.334  # N98

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 111538.22  # H107
$ws.Cells.Item(107, 9).Value = 200390  # I107
$ws.Cells.Item(107, 10).Value = 473.5  # J107
$ws.Cells.Item(107, 11).Value = 200390  # K107
$ws.Cells.Item(107, 12).Value = 473.5  # L107
$ws.Cells.Item(107, 13).Value = -198470  # M107
$ws.Cells.Item(107, 14).Value = -4313.5  # N107

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1518.6666  # H55
$ws.Cells.Item(55, 9).Value = 2657  # I55
$ws.Cells.Item(55, 10).Value = 759.7778  # J55
$ws.Cells.Item(55, 11).Value = 2657  # K55
$ws.Cells.Item(55, 12).Value = 759.7778  # L55
$ws.Cells.Item(55, 13).Value = -2484  # M55
$ws.Cells.Item(55, 14).Value = -1105.7778  # N55

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 8336366  # H61
$ws.Cells.Item(61, 9).Value = 10529016  # I61
$ws.Cells.Item(61, 10).Value = 4294  # J61
$ws.Cells.Item(61, 11).Value = 10529016  # K61
$ws.Cells.Item(61, 12).Value = 4294  # L61
$ws.Cells.Item(61, 13).Value = -10528814  # M61
$ws.Cells.Item(61, 14).Value = -4698  # N61

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(88, 8).Value = 24750  # H88
$ws.Cells.Item(88, 9).Value = 24750  # I88
$ws.Cells.Item(88, 11).Value = 24750  # K88
$ws.Cells.Item(88, 13).Value = -24322  # M88

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(91, 8).Value = 24750  # H91
$ws.Cells.Item(91, 9).Value = 24750  # I91
$ws.Cells.Item(91, 11).Value = 24750  # K91
$ws.Cells.Item(91, 13).Value = -23268  # M91

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2046.7222  # H93
$ws.Cells.Item(93, 9).Value = 1056.0667  # I93
$ws.Cells.Item(93, 11).Value = 1056.0667  # K93
$ws.Cells.Item(93, 13).Value = 191.9332999999999  # M93

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 8336366  # H113
$ws.Cells.Item(113, 9).Value = 10529016  # I113
$ws.Cells.Item(113, 10).Value = 4294  # J113
$ws.Cells.Item(113, 11).Value = 10529016  # K113
$ws.Cells.Item(113, 12).Value = 4294  # L113
$ws.Cells.Item(113, 13).Value = -10526846  # M113
$ws.Cells.Item(113, 14).Value = -8634  # N113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(114, 8).Value = 84000  # H114
$ws.Cells.Item(114, 10).Value = 84000  # J114
$ws.Cells.Item(114, 12).Value = 84000  # L114
$ws.Cells.Item(114, 14).Value = -92678  # N114

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 4578.6313  # H136
$ws.Cells.Item(136, 9).Value = 4624.375  # I136
$ws.Cells.Item(136, 11).Value = 13873.125  # K136
$ws.Cells.Item(136, 13).Value = -11323.125  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 49999  # H14
$ws.Cells.Item(14, 9).Value = 49999  # I14
$ws.Cells.Item(14, 11).Value = 49999  # K14
$ws.Cells.Item(14, 13).Value = -49831  # M14

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 454.2353  # H107
$ws.Cells.Item(107, 9).Value = 444.42856  # I107
$ws.Cells.Item(107, 11).Value = 1333.28568  # K107
$ws.Cells.Item(107, 13).Value = 586.71432  # M107

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 13965.833  # H122
$ws.Cells.Item(122, 9).Value = 15770.857  # I122
$ws.Cells.Item(122, 11).Value = 47312.571  # K122
$ws.Cells.Item(122, 13).Value = -44862.571  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3256.25  # H126
$ws.Cells.Item(126, 9).Value = 3231.7058  # I126
$ws.Cells.Item(126, 11).Value = 9695.117400000001  # K126
$ws.Cells.Item(126, 13).Value = -7225.117400000001  # M126

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 34022.87  # H132
$ws.Cells.Item(132, 9).Value = 57113.39  # I132
$ws.Cells.Item(132, 10).Value = 2051.3845  # J132
$ws.Cells.Item(132, 11).Value = 171340.17  # K132
$ws.Cells.Item(132, 12).Value = 6154.1535  # L132
$ws.Cells.Item(132, 13).Value = -168810.17  # M132
$ws.Cells.Item(132, 14).Value = -11214.1535  # N132
